$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 22:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 734525
$ws.Range("C4").Value = 24790
$ws.Range("D4").Value = 66749
$ws.Range("E4").Value = 629008
$ws.Range("G4").Value = 1614
$ws.Range("H4").Value = 38768

# Alemania (row 8)
$ws.Range("B8").Value = 143160
$ws.Range("C8").Value = 1763
$ws.Range("E8").Value = 53309
$ws.Range("G8").Value = 99
$ws.Range("H8").Value = 4451

# Canada (row 16)
$ws.Range("B16").Value = 33218
$ws.Range("C16").Value = 1291
$ws.Range("D16").Value = 11168
$ws.Range("E16").Value = 20581
$ws.Range("G16").Value = 159
$ws.Range("H16").Value = 1469

# India (row 20)
$ws.Range("B20").Value = 15722
$ws.Range("C20").Value = 1370
$ws.Range("D20").Value = 2463
$ws.Range("E20").Value = 12738
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = 521

# Sudafrica (row 53)
$ws.Range("E53").Value = 2079
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 52
